$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): insert new "Memory Usage (bytes)" column at E,
# ---- push old "Error Message" header to F ----

# E1: new header, copy header formatting (style) from D1
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Memory Usage (bytes)"

# F1: moved header, copy header formatting (style) from D1
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Error Message"

# ---- Row 2 ----
$ws.Range("D2").Value = 18.96929740905762

$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = 99031

$ws.Range("D2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# ---- Row 3 ----
$ws.Range("D3").Value = 16.99519157409668

$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = 87336

$ws.Range("D3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

# ---- Row 4 (new) ----
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2

$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = 2

$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = "Success 🟢"

$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = 20.84207534790039

$ws.Range("D3").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 87248

$ws.Range("D3").Copy()
$ws.Range("F4").PasteSpecial(-4122)

# ---- Row 5 (new) ----
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3

$ws.Range("B3").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 3

$ws.Range("C3").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "Success 🟢"

$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = 16.85214042663574

$ws.Range("D3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = 87248

$ws.Range("D3").Copy()
$ws.Range("F5").PasteSpecial(-4122)

# ---- Row 6 (new) ----
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 4

$ws.Range("B3").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = 4

$ws.Range("C3").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Success 🟢"

$ws.Range("D3").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = 16.95394515991211

$ws.Range("D3").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 87248

$ws.Range("D3").Copy()
$ws.Range("F6").PasteSpecial(-4122)

$excel.CutCopyMode = $false
